$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Calr"
$ws.Range("C2").Value = "Scarf1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 99.883077
$ws.Range("H2").Value = 299.649231
$ws.Range("I2").Value = 0.3917580942718377
$ws.Range("J2").Value = 0.3917580942718377
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 7.041396333333334
$ws.Range("N2").Value = 21.124189
$ws.Range("O2").Value = 0.5806068067430814
$ws.Range("P2").Value = 0.5806068067430814
$ws.Range("Q2").Value = 703.3163321498511
$ws.Range("R2").Value = 6329.846989348659
$ws.Range("S2").Value = 0.2274574161309267
$ws.Range("T2").Value = 0.2274574161309267

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Calr"
$ws.Range("C3").Value = "Scarf1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 99.883077
$ws.Range("H3").Value = 299.649231
$ws.Range("I3").Value = 0.3917580942718377
$ws.Range("J3").Value = 0.3917580942718377
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.539659333333333
$ws.Range("N3").Value = 4.618978
$ws.Range("O3").Value = 0.1269544628196872
$ws.Range("P3").Value = 0.1269544628196872
$ws.Range("Q3").Value = 153.785911745102
$ws.Range("R3").Value = 1384.073205705918
$ws.Range("S3").Value = 0.04973543841354552
$ws.Range("T3").Value = 0.04973543841354552

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Calr"
$ws.Range("C4").Value = "Scarf1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 99.883077
$ws.Range("H4").Value = 299.649231
$ws.Range("I4").Value = 0.3917580942718377
$ws.Range("J4").Value = 0.3917580942718377
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.546594666666666
$ws.Range("N4").Value = 10.639784
$ws.Range("O4").Value = 0.2924387304372314
$ws.Range("P4").Value = 0.2924387304372314
$ws.Range("Q4").Value = 354.244788178456
$ws.Range("R4").Value = 3188.203093606103
$ws.Range("S4").Value = 0.1145652397273654
$ws.Range("T4").Value = 0.1145652397273654

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Calr"
$ws.Range("C5").Value = "Scarf1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 124.0161413333333
$ws.Range("H5").Value = 372.048424
$ws.Range("I5").Value = 0.4864119993789693
$ws.Range("J5").Value = 0.4864119993789694
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 7.041396333333334
$ws.Range("N5").Value = 21.124189
$ws.Range("O5").Value = 0.5806068067430814
$ws.Range("P5").Value = 0.5806068067430814
$ws.Range("Q5").Value = 873.2468028586819
$ws.Range("R5").Value = 7859.221225728137
$ws.Range("S5").Value = 0.282414117720941
$ws.Range("T5").Value = 0.2824141177209411

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Calr"
$ws.Range("C6").Value = "Scarf1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 124.0161413333333
$ws.Range("H6").Value = 372.048424
$ws.Range("I6").Value = 0.4864119993789693
$ws.Range("J6").Value = 0.4864119993789694
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.539659333333333
$ws.Range("N6").Value = 4.618978
$ws.Range("O6").Value = 0.1269544628196872
$ws.Range("P6").Value = 0.1269544628196872
$ws.Range("Q6").Value = 190.9426094878525
$ws.Range("R6").Value = 1718.483485390672
$ws.Range("S6").Value = 0.06175217409020706
$ws.Range("T6").Value = 0.06175217409020707

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Calr"
$ws.Range("C7").Value = "Scarf1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 124.0161413333333
$ws.Range("H7").Value = 372.048424
$ws.Range("I7").Value = 0.4864119993789693
$ws.Range("J7").Value = 0.4864119993789694
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.546594666666666
$ws.Range("N7").Value = 10.639784
$ws.Range("O7").Value = 0.2924387304372314
$ws.Range("P7").Value = 0.2924387304372314
$ws.Range("Q7").Value = 439.8349854333795
$ws.Range("R7").Value = 3958.514868900416
$ws.Range("S7").Value = 0.1422457075678212
$ws.Range("T7").Value = 0.1422457075678212

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Calr"
$ws.Range("C8").Value = "Scarf1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 31.06188766666667
$ws.Range("H8").Value = 93.185663
$ws.Range("I8").Value = 0.121829906349193
$ws.Range("J8").Value = 0.121829906349193
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 7.041396333333334
$ws.Range("N8").Value = 21.124189
$ws.Range("O8").Value = 0.5806068067430814
$ws.Range("P8").Value = 0.5806068067430814
$ws.Range("Q8").Value = 218.7190619224786
$ws.Range("R8").Value = 1968.471557302307
$ws.Range("S8").Value = 0.07073527289121359
$ws.Range("T8").Value = 0.07073527289121362

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Calr"
$ws.Range("C9").Value = "Scarf1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 31.06188766666667
$ws.Range("H9").Value = 93.185663
$ws.Range("I9").Value = 0.121829906349193
$ws.Range("J9").Value = 0.121829906349193
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.539659333333333
$ws.Range("N9").Value = 4.618978
$ws.Range("O9").Value = 0.1269544628196872
$ws.Range("P9").Value = 0.1269544628196872
$ws.Range("Q9").Value = 47.8247252569349
$ws.Range("R9").Value = 430.422527312414
$ws.Range("S9").Value = 0.01546685031593459
$ws.Range("T9").Value = 0.01546685031593459

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Calr"
$ws.Range("C10").Value = "Scarf1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 31.06188766666667
$ws.Range("H10").Value = 93.185663
$ws.Range("I10").Value = 0.121829906349193
$ws.Range("J10").Value = 0.121829906349193
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.546594666666666
$ws.Range("N10").Value = 10.639784
$ws.Range("O10").Value = 0.2924387304372314
$ws.Range("P10").Value = 0.2924387304372314
$ws.Range("Q10").Value = 110.1639251351991
$ws.Range("R10").Value = 991.4753262167919
$ws.Range("S10").Value = 0.0356277831420448
$ws.Range("T10").Value = 0.0356277831420448
